{"js": "const replacements = [\n  [\"2026-02-06 Friday\", \"2026-02-07 Saturday\"],\n  [\"58\\u00d795=5510\", \"20\\u00d779=1580\"],\n  [\"40\\u00d768=2720\", \"63\\u00d733=2079\"],\n  [\"59\\u00d767=3953\", \"38\\u00d799=3762\"],\n  [\"16\\u00d752=832\", \"73\\u00d711=803\"],\n  [\"35\\u00d775=2625\", \"39\\u00d734=1326\"],\n  [\"62\\u00d766=4092\", \"37\\u00d741=1517\"],\n  [\"29\\u00d764=1856\", \"52\\u00d746=2392\"],\n  [\"51\\u00d714=714\", \"67\\u00d727=1809\"],\n  [\"74\\u00d762=4588\", \"16\\u00d765=1040\"],\n  [\"76\\u00d754=4104\", \"14\\u00d767=938\"],\n  [\"89\\u00d793=8277\", \"18\\u00d749=882\"],\n  [\"34\\u00d757=1938\", \"56\\u00d746=2576\"],\n  [\"70\\u00d724=1680\", \"92\\u00d741=3772\"],\n  [\"77\\u00d765=5005\", \"72\\u00d753=3816\"],\n  [\"33\\u00d725=825\", \"20\\u00d750=1000\"],\n  [\"25\\u00d796=2400\", \"25\\u00d733=825\"],\n  [\"73\\u00d712=876\", \"67\\u00d750=3350\"],\n  [\"63\\u00d757=3591\", \"62\\u00d786=5332\"],\n  [\"29\\u00d770=2030\", \"25\\u00d751=1275\"],\n  [\"68\\u00d758=3944\", \"84\\u00d777=6468\"],\n  [\"43\\u00d760=2580\", \"66\\u00d731=2046\"],\n  [\"23\\u00d753=1219\", \"93\\u00d711=1023\"],\n  [\"56\\u00d799=5544\", \"28\\u00d760=1680\"],\n  [\"70\\u00d752=3640\", \"44\\u00d713=572\"],\n  [\"44\\u00d719=836\", \"48\\u00d719=912\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for: \" + oldText);\n  }\n\n  for (const r of results.items) {\n    r.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2026-02-06 Friday\", \"2026-02-07 Saturday\"),\n    @(\"58\u00d795=5510\", \"20\u00d779=1580\"),\n    @(\"40\u00d768=2720\", \"63\u00d733=2079\"),\n    @(\"59\u00d767=3953\", \"38\u00d799=3762\"),\n    @(\"16\u00d752=832\", \"73\u00d711=803\"),\n    @(\"35\u00d775=2625\", \"39\u00d734=1326\"),\n    @(\"62\u00d766=4092\", \"37\u00d741=1517\"),\n    @(\"29\u00d764=1856\", \"52\u00d746=2392\"),\n    @(\"51\u00d714=714\", \"67\u00d727=1809\"),\n    @(\"74\u00d762=4588\", \"16\u00d765=1040\"),\n    @(\"76\u00d754=4104\", \"14\u00d767=938\"),\n    @(\"89\u00d793=8277\", \"18\u00d749=882\"),\n    @(\"34\u00d757=1938\", \"56\u00d746=2576\"),\n    @(\"70\u00d724=1680\", \"92\u00d741=3772\"),\n    @(\"77\u00d765=5005\", \"72\u00d753=3816\"),\n    @(\"33\u00d725=825\", \"20\u00d750=1000\"),\n    @(\"25\u00d796=2400\", \"25\u00d733=825\"),\n    @(\"73\u00d712=876\", \"67\u00d750=3350\"),\n    @(\"63\u00d757=3591\", \"62\u00d786=5332\"),\n    @(\"29\u00d770=2030\", \"25\u00d751=1275\"),\n    @(\"68\u00d758=3944\", \"84\u00d777=6468\"),\n    @(\"43\u00d760=2580\", \"66\u00d731=2046\"),\n    @(\"23\u00d753=1219\", \"93\u00d711=1023\"),\n    @(\"56\u00d799=5544\", \"28\u00d760=1680\"),\n    @(\"70\u00d752=3640\", \"44\u00d713=572\"),\n    @(\"44\u00d719=836\", \"48\u00d719=912\")\n)\n\nforeach ($pair in $replacements) {\n    $findText = $pair[0]\n    $replaceText = $pair[1]\n\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Text = $findText\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $replaceText\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)\n}\n"}
